$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.079.48'
$ws.Cells.Item(2, 5).Value = '  +2.06%  '
$ws.Cells.Item(3, 4).Value = '2.300.41'
$ws.Cells.Item(3, 5).Value = '  +1.40%  '
$ws.Cells.Item(4, 5).Value = '  +0.01%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '310.18'
$c.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +1.70%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '100.83'
$c.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  +4.12%  '
$ws.Cells.Item(7, 5).Value = '  +0.29%  '
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
$ws.Cells.Item(9, 5).Value = '  +3.20%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '36.16'
$c.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +1.51%  '
$ws.Cells.Item(11, 5).Value = '  +3.06%  '
$ws.Cells.Item(12, 5).Value = '  +0.73%  '
$ws.Cells.Item(13, 5).Value = '  +4.39%  '
$ws.Cells.Item(14, 4).Value = '2.661.47'
$ws.Cells.Item(14, 5).Value = '  +2.30%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '14.85'
$c.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  +3.03%  '
$ws.Cells.Item(16, 4).Value = '2.306.12'
$ws.Cells.Item(16, 5).Value = '  +1.40%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.804'
$c.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  +1.26%  '
$ws.Cells.Item(18, 4).Value = '43.094.50'
$ws.Cells.Item(18, 5).Value = '  +2.28%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '12.54'
$c.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +0.54%  '
$ws.Cells.Item(20, 5).Value = '  +0.98%  '
$ws.Cells.Item(21, 5).Value = '  +1.05%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '67.89'
$c.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +0.06%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '240.43'
$c.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +0.92%  '
$ws.Cells.Item(24, 5).Value = '  +3.82%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.61'
$c.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +1.35%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  +0.40%  '
$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '23.89'
$c.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +0.67%  '
$ws.Cells.Item(28, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '38.41'
$c.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +2.94%  '
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '2.16'
$c.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +3.02%  '
$ws.Cells.Item(30, 2).Value = 'Cosmos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '9.65'
$c.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +1.21%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '164.77'
$c.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  +3.27%  '
$ws.Cells.Item(32, 5).Value = '  +0.86%  '
$ws.Cells.Item(33, 5).Value = '  +0.16%  '
$ws.Cells.Item(34, 5).Value = '  -1.17%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '17.81'
$c.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +2.84%  '
$ws.Cells.Item(36, 5).Value = '  +0.11%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  +0.76%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.105'
$c.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -0.04%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '1.84'
$c.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +0.10%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.116'
$c.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +1.15%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '4.16'
$c.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +2.03%  '
$ws.Cells.Item(42, 5).Value = '  -5.80%  '
$ws.Cells.Item(43, 2).Value = 'VeChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.0289'
$c.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +1.59%  '
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '19.31'
$c.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +2.37%  '
$ws.Cells.Item(45, 4).Value = '1.964.12'
$ws.Cells.Item(45, 5).Value = '  -1.26%  '
$ws.Cells.Item(46, 5).Value = '  +3.18%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '9.85'
$c.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -0.86%  '
$ws.Cells.Item(48, 2).Value = 'HuobiToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '3.01'
$c.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +20.35%  '
$ws.Cells.Item(49, 2).Value = 'MultiversX'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '54.81'
$c.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  +3.13%  '
$ws.Cells.Item(50, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(50, 4).Value = '2.529.33'
$ws.Cells.Item(50, 5).Value = '  +1.60%  '
$ws.Cells.Item(51, 2).Value = 'Stacks'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '1.55'
$c.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +2.20%  '
